$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Strong on transformations and composition (Q2), trigonometric manipulation and solving (Q3), inverse functions and composite functions (Q4), geometric and arithmetic progressions (Q5b), differentiation and integration including use of limits (Q6, Q9b, Q10a), and circle geometry with tangents (Q7). Working generally clear and follows required methods."
$ws.Range("E2").Value = "Occasional algebraic/numerical accuracy slips leading to wrong final values (Q1 binomial coefficient, Q5a sum of AP, Q8b perimeter equation solved inaccurately, Q10b perpendicular bisector). Some inconsistency in deciding increasing/decreasing intervals (Q9a). Perpendicular bisector method and presentation of gradients/midpoint need more care. Ensure final answers are clearly stated and checked for reasonableness; show key intermediate values explicitly where required by the mark scheme."
$ws.Range("F2").Value = 66
$ws.Range("G2").Value = 75
$ws.Range("H2").Value = "A"

$ws.Range("D3").Value = "Strong use of calculus (differentiation/chain rule), series (AP/GP sums), integration with definite integrals and limits, and coordinate geometry. Trigonometric equation solving to final angles was correct. Good handling of inverse functions and function range concepts."
$ws.Range("E3").Value = "Transformation descriptions sometimes lacked clear order; occasional algebraic slips (e.g. Q8(b) perimeter leading to incorrect r). In Q9(a) sign/interval reasoning led to an incorrect final interval. In Q10(b) an incorrect y-value for point B propagated through the perpendicular bisector; final equation not simplified to ax+by+c=0. Careful checking of values and clarity of statements would improve accuracy."
$ws.Range("F3").Value = 64
